# Changes of 21st June 2022
# FedExShipments.xlsx — refresh the FedEx shipment test rows (rows 2-26) with
# a new batch of tracking numbers; rows that previously FAILed now PASS with
# updated ActualRate values to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> @(NewTrackingNumber, NewActualRate-or-$null, NewResult-or-$null)
$updates = @{
    2  = @("320018538422", '$19.04',  "PASS")
    3  = @("320018538433", '$27.50',  "PASS")
    4  = @("320018538466", '$31.73',  "PASS")
    5  = @("320018538488", '$43.36',  "PASS")
    6  = @("320018538525", '$56.05',  "PASS")
    7  = @("320018538547", '$231.08', "PASS")
    8  = @("320018538570", '$19.04',  "PASS")
    9  = @("320018538591", '$23.27',  "PASS")
    10 = @("320018538628", '$27.50',  "PASS")
    11 = @("320018538640", '$40.19',  "PASS")
    12 = @("320018538694", '$52.88',  "PASS")
    13 = @("320018538710", $null,     $null)
    14 = @("320018538742", $null,     $null)
    15 = @("320018538775", $null,     $null)
    16 = @("320018538801", $null,     $null)
    17 = @("320018538823", $null,     $null)
    18 = @("320018538867", $null,     $null)
    19 = @("320018538889", $null,     $null)
    20 = @("320018538915", '$62.39',  "PASS")
    21 = @("320018538937", $null,     $null)
    22 = @("320018538960", $null,     $null)
    23 = @("320018538970", $null,     $null)
    24 = @("320018538981", $null,     $null)
    25 = @("320018538992", $null,     $null)
    26 = @("320018539006", $null,     $null)
}

foreach ($row in $updates.Keys | Sort-Object) {
    $vals = $updates[$row]
    $tracking = $vals[0]
    $actualRate = $vals[1]
    $result = $vals[2]

    # Column P = ShipmentTracking. Force text so the numeric-looking
    # tracking number isn't reinterpreted as a number.
    $pCell = $ws.Cells.Item($row, 16)
    $pCell.NumberFormat = "@"
    $pCell.Value = $tracking

    if ($actualRate -ne $null) {
        # Column Q = ActualRate, stored as literal currency-formatted text.
        $qCell = $ws.Cells.Item($row, 17)
        $qCell.NumberFormat = "@"
        $qCell.Value = $actualRate
    }

    if ($result -ne $null) {
        # Column R = Result (PASS/FAIL) — plain text, no reformat needed.
        $rCell = $ws.Cells.Item($row, 18)
        $rCell.Value = $result
    }
}
